$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1/J1, copying the formatting (bold, centered, bordered)
# from the existing header cell H1 so the new headers match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-10
$values = @(
    @(3, 4),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(5, 7),
    @(5, 6),
    @(8, 8),
    @(8, 9),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
